# List of models.xlsx - fix path + fill in row 6 (R80L) data + column widths
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the previously-empty data row for model "R80L" (row 6) ---
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 15
$ws.Range("E6").Value = 415
$ws.Range("F6").Value = 50
$ws.Range("G6").Value = "DS1250"
$ws.Range("H6").Value = "Delta"
$ws.Range("I6").Value = 2850
$ws.Range("J6").Value = 0.96
$ws.Range("K6").Value = 26
$ws.Range("L6").Value = 2.5
$ws.Range("M6").Value = 4

# --- Set explicit column widths for A and B ---
$ws.Columns.Item(1).ColumnWidth = 17.5
$ws.Columns.Item(2).ColumnWidth = 12.166666666666666

# --- Update the active selection / scroll position of the sheet view ---
$ws.Range("J7").Select() | Out-Null
